$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 197.5433703333333
$ws.Range("H2").Value = 592.6301109999999
$ws.Range("I2").Value = 0.3388703761585983
$ws.Range("J2").Value = 0.3388703761585982
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.142342
$ws.Range("N2").Value = 0.427026
$ws.Range("O2").Value = 0.09118686681570291
$ws.Range("P2").Value = 0.09118686681570291
$ws.Range("Q2").Value = 28.11871841998733
$ws.Range("R2").Value = 253.068465779886
$ws.Range("S2").Value = 0.03090052785856125
$ws.Range("T2").Value = 0.03090052785856124
$ws.Range("G3").Value = 197.5433703333333
$ws.Range("H3").Value = 592.6301109999999
$ws.Range("I3").Value = 0.3388703761585983
$ws.Range("J3").Value = 0.3388703761585982
$ws.Range("M3").Value = 0.5898753333333334
$ws.Range("O3").Value = 0.377884836931721
$ws.Range("P3").Value = 0.377884836931721
$ws.Range("Q3").Value = 116.5259614231651
$ws.Range("R3").Value = 1048.733652808486
$ws.Range("S3").Value = 0.1280539768356829
$ws.Range("T3").Value = 0.1280539768356828
$ws.Range("G4").Value = 197.5433703333333
$ws.Range("H4").Value = 592.6301109999999
$ws.Range("I4").Value = 0.3388703761585983
$ws.Range("J4").Value = 0.3388703761585982
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.8287749999999999
$ws.Range("N4").Value = 2.486325
$ws.Range("O4").Value = 0.5309282962525761
$ws.Range("P4").Value = 0.5309282962525761
$ws.Range("Q4").Value = 163.7190067480083
$ws.Range("R4").Value = 1473.471060732075
$ws.Range("S4").Value = 0.1799158714643542
$ws.Range("T4").Value = 0.1799158714643541
$ws.Range("I5").Value = 0.1369374790620155
$ws.Range("J5").Value = 0.1369374790620154
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.142342
$ws.Range("N5").Value = 0.427026
$ws.Range("O5").Value = 0.09118686681570291
$ws.Range("P5").Value = 0.09118686681570291
$ws.Range("Q5").Value = 11.36277079907867
$ws.Range("R5").Value = 102.264937191708
$ws.Range("S5").Value = 0.01248689966530611
$ws.Range("T5").Value = 0.01248689966530611
$ws.Range("I6").Value = 0.1369374790620155
$ws.Range("J6").Value = 0.1369374790620154
$ws.Range("M6").Value = 0.5898753333333334
$ws.Range("O6").Value = 0.377884836931721
$ws.Range("P6").Value = 0.377884836931721
$ws.Range("Q6").Value = 47.08812727583423
$ws.Range("S6").Value = 0.05174659694519068
$ws.Range("T6").Value = 0.05174659694519067
$ws.Range("I7").Value = 0.1369374790620155
$ws.Range("J7").Value = 0.1369374790620154
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.8287749999999999
$ws.Range("N7").Value = 2.486325
$ws.Range("O7").Value = 0.5309282962525761
$ws.Range("P7").Value = 0.5309282962525761
$ws.Range("Q7").Value = 66.15883132881666
$ws.Range("R7").Value = 595.4294819593499
$ws.Range("S7").Value = 0.07270398245151868
$ws.Range("T7").Value = 0.07270398245151867
$ws.Range("G8").Value = 148.824417
$ws.Range("H8").Value = 446.473251
$ws.Range("I8").Value = 0.2552967790580629
$ws.Range("J8").Value = 0.2552967790580629
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.142342
$ws.Range("N8").Value = 0.427026
$ws.Range("O8").Value = 0.09118686681570291
$ws.Range("P8").Value = 0.09118686681570291
$ws.Range("Q8").Value = 21.183965164614
$ws.Range("R8").Value = 190.655686481526
$ws.Range("S8").Value = 0.02327971339044551
$ws.Range("T8").Value = 0.02327971339044551
$ws.Range("G9").Value = 148.824417
$ws.Range("H9").Value = 446.473251
$ws.Range("I9").Value = 0.2552967790580629
$ws.Range("J9").Value = 0.2552967790580629
$ws.Range("M9").Value = 0.5898753333333334
$ws.Range("O9").Value = 0.377884836931721
$ws.Range("P9").Value = 0.377884836931721
$ws.Range("Q9").Value = 87.78785258601403
$ws.Range("R9").Value = 790.0906732741261
$ws.Range("S9").Value = 0.0964727817235497
$ws.Range("T9").Value = 0.0964727817235497
$ws.Range("G10").Value = 148.824417
$ws.Range("H10").Value = 446.473251
$ws.Range("I10").Value = 0.2552967790580629
$ws.Range("J10").Value = 0.2552967790580629
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.8287749999999999
$ws.Range("N10").Value = 2.486325
$ws.Range("O10").Value = 0.5309282962525761
$ws.Range("P10").Value = 0.5309282962525761
$ws.Range("Q10").Value = 123.341956199175
$ws.Range("R10").Value = 1110.077605792575
$ws.Range("S10").Value = 0.1355442839440677
$ws.Range("T10").Value = 0.1355442839440677
$ws.Range("G11").Value = 35.426853
$ws.Range("H11").Value = 106.280559
$ws.Range("I11").Value = 0.06077202683121193
$ws.Range("J11").Value = 0.06077202683121192
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.142342
$ws.Range("N11").Value = 0.427026
$ws.Range("O11").Value = 0.09118686681570291
$ws.Range("P11").Value = 0.09118686681570291
$ws.Range("Q11").Value = 5.042729109726
$ws.Range("R11").Value = 45.384561987534
$ws.Range("S11").Value = 0.005541610716778046
$ws.Range("T11").Value = 0.005541610716778046
$ws.Range("G12").Value = 35.426853
$ws.Range("H12").Value = 106.280559
$ws.Range("I12").Value = 0.06077202683121193
$ws.Range("J12").Value = 0.06077202683121192
$ws.Range("M12").Value = 0.5898753333333334
$ws.Range("O12").Value = 0.377884836931721
$ws.Range("P12").Value = 0.377884836931721
$ws.Range("Q12").Value = 20.897426722326
$ws.Range("R12").Value = 188.076840500934
$ws.Range("S12").Value = 0.02296482744912269
$ws.Range("T12").Value = 0.02296482744912269
$ws.Range("G13").Value = 35.426853
$ws.Range("H13").Value = 106.280559
$ws.Range("I13").Value = 0.06077202683121193
$ws.Range("J13").Value = 0.06077202683121192
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.8287749999999999
$ws.Range("N13").Value = 2.486325
$ws.Range("O13").Value = 0.5309282962525761
$ws.Range("P13").Value = 0.5309282962525761
$ws.Range("Q13").Value = 29.360890095075
$ws.Range("R13").Value = 264.248010855675
$ws.Range("S13").Value = 0.03226558866531119
$ws.Range("T13").Value = 0.03226558866531119
$ws.Range("G14").Value = 121.3248153333333
$ws.Range("H14").Value = 363.974446
$ws.Range("I14").Value = 0.2081233388901116
$ws.Range("J14").Value = 0.2081233388901115
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.142342
$ws.Range("N14").Value = 0.427026
$ws.Range("O14").Value = 0.09118686681570291
$ws.Range("P14").Value = 0.09118686681570291
$ws.Range("Q14").Value = 17.26961686417733
$ws.Range("R14").Value = 155.426551777596
$ws.Range("S14").Value = 0.01897811518461201
$ws.Range("T14").Value = 0.018978115184612
$ws.Range("G15").Value = 121.3248153333333
$ws.Range("H15").Value = 363.974446
$ws.Range("I15").Value = 0.2081233388901116
$ws.Range("J15").Value = 0.2081233388901115
$ws.Range("M15").Value = 0.5898753333333334
$ws.Range("O15").Value = 0.377884836931721
$ws.Range("P15").Value = 0.377884836931721
$ws.Range("Q15").Value = 71.56651588635512
$ws.Range("R15").Value = 644.0986429771961
$ws.Range("S15").Value = 0.07864665397817512
$ws.Range("T15").Value = 0.07864665397817511
$ws.Range("G16").Value = 121.3248153333333
$ws.Range("H16").Value = 363.974446
$ws.Range("I16").Value = 0.2081233388901116
$ws.Range("J16").Value = 0.2081233388901115
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.8287749999999999
$ws.Range("N16").Value = 2.486325
$ws.Range("O16").Value = 0.5309282962525761
$ws.Range("P16").Value = 0.5309282962525761
$ws.Range("Q16").Value = 123.341956199175
$ws.Range("R16").Value = 1110.077605792575
$ws.Range("S16").Value = 0.1355442839440677
$ws.Range("T16").Value = 0.1355442839440677
